$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$hdr = $ws.Range("A1:C1")

# Mark the header row as locked (protection) and force text storage so the
# numeric-looking "4" in A1 is kept as text, matching the rest of the row.
$hdr.NumberFormat = "@"
$hdr.Locked = $true

# Replace the header labels (id / customerId / itemNotFound) with the new
# record values.
$ws.Range("A1").Value = "4"
$ws.Range("B1").Value = "Scooby Snacks"
$ws.Range("C1").Value = ""
